$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3835.182
$ws.Range("J43").Value = 3828
$ws.Range("L43").Value = 3828
$ws.Range("N43").Value = -3966

$ws.Range("H86").Value = 6210.8423
$ws.Range("I86").Value = 7844.1816
$ws.Range("J86").Value = 3965
$ws.Range("K86").Value = 7844.1816
$ws.Range("L86").Value = 3965
$ws.Range("M86").Value = -6721.1816
$ws.Range("N86").Value = -6211

$ws.Range("H89").Value = 6210.8423
$ws.Range("I89").Value = 7844.1816
$ws.Range("J89").Value = 3965
$ws.Range("K89").Value = 39220.908
$ws.Range("L89").Value = 19825
$ws.Range("M89").Value = -33604.908
$ws.Range("N89").Value = -31057

$ws.Range("H98").Value = 1189.92
$ws.Range("I98").Value = 1189.92
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1189.92
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = 308.0799999999999

$ws.Range("H122").Value = 1189.92
$ws.Range("I122").Value = 1189.92
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3569.76
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -1119.76

$ws.Range("H129").Value = 2999.7778
$ws.Range("J129").Value = 3124.75
$ws.Range("L129").Value = 9374.25
$ws.Range("N129").Value = -19374.25

$ws.Range("H132").Value = 8510.267
$ws.Range("I132").Value = 9021.75
$ws.Range("J132").Value = 1349.5
$ws.Range("K132").Value = 27065.25
$ws.Range("L132").Value = 4048.5
$ws.Range("M132").Value = -24535.25
$ws.Range("N132").Value = -9108.5

$ws.Range("H135").Value = 12000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 12000
$ws.Range("K135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("M135").Value = 108000
$ws.Range("N135").Value = -113070

$ws.Range("H137").Value = 5891273.5
$ws.Range("I137").Value = 11113879
$ws.Range("K137").Value = 33341637
$ws.Range("M137").Value = -33339087

$ws.Range("H138").Value = 5927.8906
$ws.Range("I138").Value = 5906.1333
$ws.Range("J138").Value = 5934.551
$ws.Range("K138").Value = 17718.3999
$ws.Range("L138").Value = 17803.653
$ws.Range("M138").Value = -12578.3999
$ws.Range("N138").Value = -28083.653

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3271.8157
$ws.Range("I132").Value = 2131.7693
$ws.Range("K132").Value = 6395.3079
$ws.Range("M132").Value = -3865.3079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1698.9131
$ws.Range("I86").Value = 1539.6
$ws.Range("K86").Value = 1539.6
$ws.Range("M86").Value = -416.5999999999999

$ws.Range("H89").Value = 1698.9131
$ws.Range("I89").Value = 1539.6
$ws.Range("K89").Value = 7698
$ws.Range("M89").Value = -2082

$ws.Range("H134").Value = 6169145
$ws.Range("I134").Value = 5007228
$ws.Range("J134").Value = 9267591
$ws.Range("K134").Value = 15021684
$ws.Range("L134").Value = 27802773
$ws.Range("M134").Value = -15019149
$ws.Range("N134").Value = -27807843

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").ClearContents()
$ws.Range("N20").Value = 0

$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").ClearContents()
$ws.Range("N30").Value = 0

$ws.Range("H31").Value = 364186.1
$ws.Range("I31").Value = 911164.4399999999
$ws.Range("J31").Value = 3677.6365
$ws.Range("K31").Value = 911164.4399999999
$ws.Range("L31").Value = 3677.6365
$ws.Range("M31").Value = -910869.4399999999
$ws.Range("N31").Value = -4267.636500000001

$ws.Range("H34").Value = 364186.1
$ws.Range("I34").Value = 911164.4399999999
$ws.Range("J34").Value = 3677.6365
$ws.Range("K34").Value = 911164.4399999999
$ws.Range("L34").Value = 3677.6365
$ws.Range("M34").Value = -910962.4399999999
$ws.Range("N34").Value = -4081.6365

$ws.Range("H99").Value = 754657.25
$ws.Range("I99").Value = 38678.875
$ws.Range("K99").Value = 38678.875
$ws.Range("M99").Value = -37180.875

$ws.Range("H122").Value = 29507.555
$ws.Range("J122").Value = 52389.8
$ws.Range("L122").Value = 157169.4
$ws.Range("N122").Value = -162069.4

$ws.Range("H123").Value = 128000
$ws.Range("J123").Value = 128000
$ws.Range("L123").Value = 128000
$ws.Range("N123").Value = -137800

$ws.Range("H126").Value = 754657.25
$ws.Range("I126").Value = 38678.875
$ws.Range("K126").Value = 116036.625
$ws.Range("M126").Value = -113566.625

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").ClearContents()
$ws.Range("N128").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 6851.2856
$ws.Range("J39").Value = 11791.25
$ws.Range("L39").Value = 35373.75
$ws.Range("N39").Value = -35961.75

$ws.Range("H70").Value = 14888
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 14888
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H75").Value = 995
$ws.Range("J75").Value = 995
$ws.Range("L75").Value = 2985
$ws.Range("N75").Value = -4981

$ws.Range("H78").Value = 995
$ws.Range("J78").Value = 995
$ws.Range("L78").Value = 8955
$ws.Range("N78").Value = -18939

$ws.Range("H88").Value = 40550.668
$ws.Range("I88").Value = 54826
$ws.Range("J88").Value = 12000
$ws.Range("K88").Value = 164478
$ws.Range("L88").Value = 36000
$ws.Range("M88").Value = -164050
$ws.Range("N88").Value = -36856

$ws.Range("H91").Value = 40550.668
$ws.Range("I91").Value = 54826
$ws.Range("J91").Value = 12000
$ws.Range("K91").Value = 164478
$ws.Range("L91").Value = 36000
$ws.Range("M91").Value = -162996
$ws.Range("N91").Value = -38964

$ws.Range("H97").Value = 712.63635
$ws.Range("I97").Value = 670
$ws.Range("J97").Value = 722.1111
$ws.Range("K97").Value = 2010
$ws.Range("L97").Value = 2166.3333
$ws.Range("M97").Value = -1514
$ws.Range("N97").Value = -3158.3333

$ws.Range("H132").Value = 769.25
$ws.Range("J132").Value = 778.4286
$ws.Range("L132").Value = 7005.8574
$ws.Range("N132").Value = -12065.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1826.25
$ws.Range("I80").Value = 2000
$ws.Range("K80").Value = 2000
$ws.Range("M80").Value = -1002

$ws.Range("H83").Value = 1826.25
$ws.Range("I83").Value = 2000
$ws.Range("K83").Value = 10000
$ws.Range("M83").Value = -5008

$ws.Range("H102").Value = 3065.875
$ws.Range("I102").Value = 1408.5714
$ws.Range("K102").Value = 1408.5714
$ws.Range("M102").Value = 213.4286

$ws.Range("H107").Value = 1067.7587
$ws.Range("J107").Value = 1133.5714
$ws.Range("L107").Value = 1133.5714
$ws.Range("N107").Value = -4973.5714

$ws.Range("H113").Value = 1468.8889
$ws.Range("I113").Value = 1561.0625
$ws.Range("J113").Value = 731.5
$ws.Range("K113").Value = 1561.0625
$ws.Range("L113").Value = 731.5
$ws.Range("M113").Value = 608.9375
$ws.Range("N113").Value = -5071.5

$ws.Range("H122").Value = 2842384
$ws.Range("I122").Value = 3552229.8
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 10656689.4
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -10654239.4
$ws.Range("N122").Value = -13900

$ws.Range("H141").Value = 95999.664
$ws.Range("J141").Value = 95999.664
$ws.Range("L141").Value = 95999.664
$ws.Range("N141").Value = -106359.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 669
$ws.Range("I7").Value = 669
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 669
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -557

$ws.Range("H16").Value = 1599
$ws.Range("I16").Value = 1599
$ws.Range("J16").Value = 1599
$ws.Range("K16").Value = 1599
$ws.Range("L16").Value = 1599
$ws.Range("M16").Value = -1429
$ws.Range("N16").Value = -1939

$ws.Range("H40").Value = 6804
$ws.Range("I40").Value = 6405.3335
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 6405.3335
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -6269.3335
$ws.Range("N40").Value = -8272

$ws.Range("H126").Value = 669
$ws.Range("I126").Value = 669
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2007
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = 463

$ws.Range("H131").Value = 57714.5
$ws.Range("J131").Value = 57714.5
$ws.Range("L131").Value = 57714.5
$ws.Range("N131").Value = -67794.5

$ws.Range("H132").Value = 2877263
$ws.Range("I132").Value = 4169976.5
$ws.Range("K132").Value = 12509929.5
$ws.Range("M132").Value = -12507399.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 46823.84
$ws.Range("I122").Value = 1623.9375
$ws.Range("J122").Value = 127179.22
$ws.Range("K122").Value = 4871.8125
$ws.Range("L122").Value = 381537.66
$ws.Range("M122").Value = -2421.8125
$ws.Range("N122").Value = -386437.66

$ws.Range("H127").Value = 69935
$ws.Range("J127").Value = 69935
$ws.Range("L127").Value = 69935
$ws.Range("N127").Value = -79855

$ws.Range("H136").Value = 2525693.5
$ws.Range("I136").Value = 1673813.6
$ws.Range("K136").Value = 5021440.800000001
$ws.Range("M136").Value = -5018890.800000001
